$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: update the generated date/time stamp under the title.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("November   3, 2021 (06:51:59 PM)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "November   3, 2021 (11:34:30 PM)", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: reword the "Problem 6" prompt and split the sentence so that the
# sentinel value "Done" is called out in its own run (surrounded by curly
# quotes), matching the target run layout:
#
#   [0] "Ask the user to enter integers. ... they are done (by entering a
#        sentinel value like"
#   [1] " "
#   [2] "\u201c"
#   [3] "Done"
#   [4] "\u201d"
#   [5] "), display the smallest value the user entered. If the user did not
#        enter any integers, display"
#   [6..9] (unchanged) " ", "\u201c", "You did not enter anything.", "\u201d"
# ---------------------------------------------------------------------------

# Locate the paragraph that starts with the old sentence so we do not have to
# hard-code paragraph indices.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Ask user to enter integers*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne -1) {
    $oldFirstRun = "Ask user to enter integers. Keep track of the smallest value user enters. After user indicates they are done, display the smallest value user entered. If user did not enter any integers display"
    $newFirstRun = "Ask the user to enter integers. Keep track of the smallest value the user enters. After the user indicates they are done (by entering a sentinel value like"

    $para = $d.Paragraphs($targetIndex)
    $paraStart = $para.Range.Start

    # Split the paragraph right after the first run (i.e. right before the
    # " " run that precedes the opening curly quote). This isolates the run
    # we need to rewrite into its own paragraph without disturbing the runs
    # that must stay untouched.
    $splitPos = $paraStart + $oldFirstRun.Length
    $d.Range($splitPos, $splitPos).InsertParagraphAfter() | Out-Null

    # Rewrite the isolated first run's text.
    $paraA = $d.Paragraphs($targetIndex)
    $rngA = $d.Range($paraA.Range.Start, $paraA.Range.End - 1)
    $rngA.Text = $newFirstRun

    # Append the remaining new pieces; because paraA only ever contains one
    # run at this point, each InsertAfter merges into that same run, which
    # lets us build up the full combined text we will split into separate
    # runs in the next step.
    $pieces = @(" ", [char]0x201C, "Done", [char]0x201D, "), display the smallest value the user entered. If the user did not enter any integers, display")
    foreach ($piece in $pieces) {
        $cur = $d.Paragraphs($targetIndex)
        $endPos = $cur.Range.End - 1
        $d.Range($endPos, $endPos).InsertAfter($piece) | Out-Null
    }

    # Compute the relative offsets (within paraA) at which the six runs
    # should be split apart, then split from the highest offset down to the
    # lowest so earlier offsets stay valid.
    $newRunTexts = @($newFirstRun, " ", [char]0x201C, "Done", [char]0x201D, "), display the smallest value the user entered. If the user did not enter any integers, display")
    $offsets = New-Object System.Collections.ArrayList
    $pos = 0
    for ($k = 0; $k -lt $newRunTexts.Length - 1; $k++) {
        $pos += $newRunTexts[$k].Length
        [void]$offsets.Add($pos)
    }

    $paraFinal = $d.Paragraphs($targetIndex)
    $baseStart = $paraFinal.Range.Start

    for ($k = $offsets.Count - 1; $k -ge 0; $k--) {
        $p = $baseStart + $offsets[$k]
        $d.Range($p, $p).InsertParagraphAfter() | Out-Null
    }

    # Merge the newly created paragraphs (and the following, untouched
    # paragraph that still holds the original trailing runs) back into a
    # single paragraph by deleting the paragraph marks between them. This
    # keeps every run distinct instead of collapsing them into one. There
    # are (newRunTexts.Length - 1) marks between the six new-run paragraphs
    # plus one more mark joining the last of them to the paragraph that
    # still holds the original, untouched trailing runs.
    $mergeCount = $newRunTexts.Length
    for ($k = 0; $k -lt $mergeCount; $k++) {
        $pp = $d.Paragraphs($targetIndex)
        $pmStart = $pp.Range.End - 1
        $pmEnd = $pp.Range.End
        $d.Range($pmStart, $pmEnd).Delete() | Out-Null
    }
}
